# Updated symbol list refresh (Sun Dec 25 22:54:05 UTC 2022, GitHub Actions).
# The source sheet stores every "Price"/"Volume(1h)" figure as literal text
# (inline strings), so numeric-looking values are written with a leading
# apostrophe to keep them as text instead of being parsed as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'242.88"
$ws.Range("D3").Value = "'23.08"
$ws.Range("D4").Value = "'5.393"
$ws.Range("D6").Value = "'3.400"
$ws.Range("D7").Value = "'6.477"
$ws.Range("D8").Value = "'0.8130"
$ws.Range("D9").Value = "'0.9081"
$ws.Range("D10").Value = "'0.1415"
$ws.Range("D11").Value = "'0.07426"
$ws.Range("D12").Value = "'0.03330"
$ws.Range("D13").Value = "'0.03064"
$ws.Range("D14").Value = "'0.09334"
$ws.Range("D15").Value = "'3.863"
$ws.Range("D16").Value = "'0.001584"
$ws.Range("D17").Value = "'0.04633"
$ws.Range("D18").Value = "'0.0005938"
$ws.Range("D19").Value = "'0.006107"
$ws.Range("D20").Value = "'0.005028"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("D21").Value = "'0.0009868"
$ws.Range("D23").Value = "'0.0002899"
$ws.Range("D24").Value = "'3.625"
$ws.Range("D25").Value = "'2.162"
$ws.Range("D27").Value = "'0.1293"
$ws.Range("D40").Value = "'0.03883"
$ws.Range("D41").Value = "'0.006208"
$ws.Range("D42").Value = "'0.1072"
$ws.Range("D43").Value = "'0.002799"
$ws.Range("D44").Value = "'0.007216"
$ws.Range("D45").Value = "'0.00005198"
$ws.Range("D47").Value = "'0.0005798"
$ws.Range("E48").Value = "47CoinbaseStockTokenCOINBestin24h"
$ws.Range("D49").Value = "'0.002260"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("D51").Value = "'0.0001999"
